# Append new subscriber rows (34-40) to the worksheet, matching the
# six existing columns (first name, last name, national ID, phone,
# region, specialty).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ERE",  "FS",   "14000000", "98625905", "GRG",  "السقوي"),
    @("RR",   "RR",   "12345678", "22334455", "EE",   "الصيد بالأضواء"),
    @("jjj",  "jjj",  "12345678", "12345678", "TGFF", "الماشية"),
    @("UUUU", "YYYY", "12335511", "11223344", "HHH",  "الصيد الساحلي"),
    @("ggf",  "ff",   "12345678", "32323323", "222",  "السقوي"),
    @("RRR",  "CGG",  "12345678", "12345678", "rrrr", "الزياتين"),
    @("433",  "DDD",  "12121212", "12121212", "FFF",  "الصيد بالأضواء")
)

$startRow = 34
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $val = $rowData[$col - 1]
        if ($val -match '^[0-9]+$') {
            # Purely-numeric values (IDs, phone numbers, ...) must stay
            # text, like every other such value already in the sheet -
            # force Excel to keep the leading apostrophe / text origin.
            $cell.Value = "'" + $val
        } else {
            $cell.Value = $val
        }
    }
}
